# Vinayak: Change to add all items in cart.
# Adds a new test case row (GK_016_Test) to the TestCases sheet and
# leaves the workbook focused on that sheet/row, matching the authored edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCases")

# New test case: row 21
$ws.Range("A21").Value = "GK_016_Test"
$ws.Range("B21").Value = "Verify user is able to add all vegitables."

# Bring TestCases to the front as the active sheet and reflect the
# author's final selection/scroll position on it.
$ws.Activate() | Out-Null
$ws.Range("C20").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1

$wb.Save() | Out-Null
